$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.329.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.925.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "468.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.731"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000338"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.561.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.922.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.568.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "38.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.49%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.43%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "726.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +17.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0478"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("E42").Value = "  -5.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("E47").Value = "  +5.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.07%  "
